$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 summary cells ---
$ws.Range("D2").Value = "8.333"
$ws.Range("F2").Value = "analysis : 1"
$ws.Range("G2").Value = "9.090"
$ws.Range("I2").Value = "0.0"
$ws.Range("K2").Value = "0.0"

# --- Remove stray cluster-match cells in F3:F5 ---
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()

# --- Rewrite the CV KEYWORDS column (E2:E57) with the refreshed keyword counts ---
$keywords = @(
  "data engineer : 1",
  "engineer : 1",
  "apache : 6",
  "celery : 1",
  "aws : 16",
  "lambda : 1",
  "research : 1",
  "data migration : 1",
  "migration : 1",
  "mongodb : 3",
  "s3 : 1",
  "kubernetes : 3",
  "kinesis : 1",
  "dynamodb : 3",
  "reports : 1",
  "qlik : 4",
  "conversion : 1",
  "sql : 1",
  "ssis : 1",
  "etl : 1",
  "selenium : 2",
  "python : 5",
  "sql queries : 1",
  "queries : 1",
  "stored procedures : 1",
  "bi : 2",
  "developer : 3",
  "bamboo : 1",
  "deployment : 1",
  "software developer : 3",
  "analysis : 3",
  "data mining : 3",
  "mining : 2",
  "tensorflow : 3",
  "intern : 1",
  "computer engineering : 1",
  "electrical : 1",
  "electronics : 1",
  "hybrid : 1",
  "design : 1",
  "tools : 1",
  "azure : 2",
  "amazon : 1",
  "microsoft azure : 1",
  "unix : 1",
  "java : 1",
  "tableau : 1",
  "english : 1",
  "apache spark : 2",
  "github : 1",
  "docker : 1",
  "ansible : 1",
  "amazon web services : 1",
  "redshift : 1",
  "mysql : 1",
  "unix shell : 1"
)

for ($i = 0; $i -lt $keywords.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 5).Value = $keywords[$i]
}

# --- Drop the now-unused stop-word rows at the bottom (old rows 58:63) ---
$ws.Range("A58:K63").EntireRow.Delete() | Out-Null
